$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.030.92'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').Value = '3.723.86'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.89'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +8.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '193.12'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +11.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.640'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.726'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '60.32'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +18.29%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.47'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '4.314.36'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '3.725.58'
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.16'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.64%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.52'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.99'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').Value = '68.886.58'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '413.22'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('E22').Value = '  +3.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '90.30'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.12'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.20'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.47%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.39'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +7.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.84'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.06'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.36%  '
$ws.Range('E29').Value = '  +4.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.99'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.80'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.79'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.21%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '46.84'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +9.33%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.124'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +7.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '640.58'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +11.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '67.64'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.50%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0835'
$ws.Range('E37').Value = '  -6.51%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.417'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +6.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('E41').Value = '  +6.29%  '
$ws.Range('E42').Value = '  +3.99%  '
$ws.Range('E43').Value = '  +3.81%  '
$ws.Range('E44').Value = '  +3.52%  '
$ws.Range('D45').Value = '2.932.09'
$ws.Range('E45').Value = '  +8.92%  '
$ws.Range('E46').Value = '  +5.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('E48').Value = '  +2.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '145.00'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.56'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -13.29%  '
